$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Thbs2"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.8911683333333333
$ws.Range("H2").Value = 2.673505
$ws.Range("I2").Value = 0.02693425114262819
$ws.Range("J2").Value = 0.02693425114262819
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 100.26837202632
$ws.Range("R2").Value = 902.4153482368799
$ws.Range("S2").Value = 0.008821810230315415
$ws.Range("T2").Value = 0.008821810230315415
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Thbs2"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8911683333333333
$ws.Range("H3").Value = 2.673505
$ws.Range("I3").Value = 0.02693425114262819
$ws.Range("J3").Value = 0.02693425114262819
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 94.74408547444332
$ws.Range("R3").Value = 852.69676926999
$ws.Range("S3").Value = 0.008335772543319291
$ws.Range("T3").Value = 0.008335772543319292
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Thbs2"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8911683333333333
$ws.Range("H4").Value = 2.673505
$ws.Range("I4").Value = 0.02693425114262819
$ws.Range("J4").Value = 0.02693425114262819
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 111.1212546639811
$ws.Range("R4").Value = 1000.09129197583
$ws.Range("S4").Value = 0.009776668368993477
$ws.Range("T4").Value = 0.009776668368993479
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Thbs2"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 24.359699
$ws.Range("H5").Value = 73.07909699999999
$ws.Range("I5").Value = 0.7362360466408275
$ws.Range("J5").Value = 0.7362360466408276
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 2740.792362589008
$ws.Range("R5").Value = 24667.13126330107
$ws.Range("S5").Value = 0.2411403477969229
$ws.Range("T5").Value = 0.241140347796923
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Thbs2"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 24.359699
$ws.Range("H6").Value = 73.07909699999999
$ws.Range("I6").Value = 0.7362360466408275
$ws.Range("J6").Value = 0.7362360466408276
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 2589.788391105733
$ws.Range("R6").Value = 23308.0955199516
$ws.Range("S6").Value = 0.2278547189038985
$ws.Range("T6").Value = 0.2278547189038985
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Thbs2"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 24.359699
$ws.Range("H7").Value = 73.07909699999999
$ws.Range("I7").Value = 0.7362360466408275
$ws.Range("J7").Value = 0.7362360466408276
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 3037.451191731745
$ws.Range("R7").Value = 27337.0607255857
$ws.Range("S7").Value = 0.2672409799400061
$ws.Range("T7").Value = 0.2672409799400062
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Thbs2"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.835938333333334
$ws.Range("H8").Value = 23.507815
$ws.Range("I8").Value = 0.2368297022165442
$ws.Range("J8").Value = 0.2368297022165442
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 881.64800138616
$ws.Range("R8").Value = 7934.83201247544
$ws.Range("S8").Value = 0.07756913970961797
$ws.Range("T8").Value = 0.07756913970961797
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Thbs2"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.835938333333334
$ws.Range("H9").Value = 23.507815
$ws.Range("I9").Value = 0.2368297022165442
$ws.Range("J9").Value = 0.2368297022165442
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 833.0735995172633
$ws.Range("R9").Value = 7497.66239565537
$ws.Range("S9").Value = 0.07329546749694854
$ws.Range("T9").Value = 0.07329546749694855
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Thbs2"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 7.835938333333334
$ws.Range("H10").Value = 23.507815
$ws.Range("I10").Value = 0.2368297022165442
$ws.Range("J10").Value = 0.2368297022165442
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 977.0761218732546
$ws.Range("R10").Value = 8793.68509685929
$ws.Range("S10").Value = 0.0859650950099777
$ws.Range("T10").Value = 0.08596509500997772
